$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style swap: three tables (slides 14, 15, 16) move from the local
#    "Table_0" style to the built-in "{844D684E-9CB8-4158-B53E-804B2317C908}"
#    style, matching <a:tableStyleId> in each table's <a:tblPr>.
# ---------------------------------------------------------------------------
$newStyleId = "{844D684E-9CB8-4158-B53E-804B2317C908}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Theme colour swap: the deck's live theme (the one bound to the slide
#    master / presentation, serialized as ppt/theme/theme2.xml) changes its
#    palette from the "Red Violet" (Integral) scheme to the "Office" scheme
#    that used to live in ppt/theme/theme1.xml.
# ---------------------------------------------------------------------------
function ToOleColor([int]$r, [int]$g, [int]$b) {
    return $b * 65536 + $g * 256 + $r
}

$officePalette = @(
    @(0x00, 0x00, 0x00),  # dk1
    @(0xFF, 0xFF, 0xFF),  # lt1
    @(0x44, 0x54, 0x6A),  # dk2
    @(0xE7, 0xE6, 0xE6),  # lt2
    @(0x5B, 0x9B, 0xD5),  # accent1
    @(0xED, 0x7D, 0x31),  # accent2
    @(0xA5, 0xA5, 0xA5),  # accent3
    @(0xFF, 0xC0, 0x00),  # accent4
    @(0x44, 0x72, 0xC4),  # accent5
    @(0x70, 0xAD, 0x47),  # accent6
    @(0x05, 0x63, 0xC1),  # hlink
    @(0x95, 0x4F, 0x72)   # folHlink
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le $officePalette.Count; $i++) {
    $rgb = $officePalette[$i - 1]
    $colorScheme.Colors($i).RGB = ToOleColor $rgb[0] $rgb[1] $rgb[2]
}
